$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the placeholder-styled empty cells before re-populating
$ws.Range("A3").Clear()
$ws.Range("A4").Clear()
$ws.Range("A5").Clear()

# 1) New row 4 - Keystone Electronics PCB terminals
$ws.Range("A16").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = "1287-ST"

$ws.Range("A16").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("B4").Value = "(534-1287-ST)"

# 2) Header row: "Digikey #" -> "Digikey # (Mouser)"
$ws.Range("B1").Value = "Digikey # (Mouser)"

$ws.Range("C4").Value = "Terminals PCB STURDI-MNT TERM"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = "Keystone Electronics"
$ws.Range("G4").Value = "Mouser"
$ws.Range("I4").Value = "https://www.mouser.com/ProductDetail/Keystone-Electronics/1287-ST?qs=lQmX4aIt5iBQmFQ9gmrtHw%3D%3D"

$ws.Rows.Item(4).RowHeight = 15.75

# 3) New row 3 - TE Connectivity cut-strip terminals
$ws.Range("A16").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A3").Value = "42100-2 (CUT STRIP)"

# B3 uses a new style: same font as A3/A16 plus left/center/wrap/indent alignment
$ws.Range("A16").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("B3").HorizontalAlignment = -4131
$ws.Range("B3").VerticalAlignment = -4108
$ws.Range("B3").WrapText = $true
$ws.Range("B3").IndentLevel = 1
$ws.Range("B3").Value = "(571-421002-CT)"

$ws.Range("C3").Value = "Terminals .25 FF REC IS 18-14 Cut Strip of 100"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = "TE Connectivity"
$ws.Range("G3").Value = "Mouser"
$ws.Range("I3").Value = "https://www.mouser.com/ProductDetail/TE-Connectivity-AMP/42100-2-CUT-STRIP?qs=2FIyTMJ0hNlByrnrD71s6A%3D%3D"

$ws.Rows.Item(3).RowHeight = 15.75

# 4) Update the active selection shown when the sheet is opened
$ws.Range("H6").Select()
